# Migracion de casos UFT
# - Corrige Sucursal del usuario F00743 (fila 12) de "037" a "043"
# - Agrega 5 usuarios nuevos para los nuevos modulos de prueba (APT, CHR, COM, CICC, CDC)
# - Deja la hoja "Users" como hoja activa (en vez de "Carga Saldos")

$wb = $excel.ActiveWorkbook
$wsUsers = $wb.Worksheets.Item("Users")

# --- Corrige el valor de Sucursal en C12 (era "037", ahora "043") ---
$wsUsers.Range("C12").Value = "043"

# --- Nuevas filas de usuarios (Administracion de Piezas, Cheques Rechazados,
#     Comisiones, Consulta Chq Ingresados por Camara y Canje, Cierre de cuenta) ---
$newUsers = @(
    @{ Row = 44; User = "F00068";   Sucursal = "068" },
    @{ Row = 45; User = "MSORACE";  Sucursal = "089" },
    @{ Row = 46; User = "F00419";   Sucursal = "019" },
    @{ Row = 47; User = "ATORRA";   Sucursal = "Casa central" },
    @{ Row = 48; User = "CRECERAB"; Sucursal = "Usuario Emergencia" }
)

foreach ($u in $newUsers) {
    $r = $u.Row
    $wsUsers.Range("A$r").Value = $u.User

    $cCell = $wsUsers.Range("C$r")
    $cCell.NumberFormat = "@"
    $cCell.HorizontalAlignment = -4152
    $cCell.Value = $u.Sucursal
}

# --- Deja "Users" como hoja/selección activa del libro ---
$wsUsers.Activate()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$wsUsers.Range("D16").Select()
